# fix protocol case code
# Row 18 (protocol_017): parameterize the variable name in the title, and
# make the set/show statements operate on the GLOBAL variable.
# Row 19 (protocol_018): disable the case (Testable = n) and parameterize
# the variable name in the title too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C18").Value = "设置变量值interactive_timeout"
$ws.Range("H18").Value = "set global interactive_timeout=14400"
$ws.Range("I18").Value = "show global variables like 'interactive_timeout'"

$ws.Range("B19").Value = "n"
$ws.Range("C19").Value = "设置变量值wait_timeout"

# Restore the view: selection moved to H24 and the zoom level changed to 100%.
[void]$ws.Activate()
$ws.Range("H24").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100
